# Edit script: apply the changes described by the diff to final_invoice_data.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the old "Cost_Type" (R) and "Remarks" (S) columns. ---
# Deleting entire columns shifts everything to the right of them left,
# turning the old layout (... P Q R S T) into (... P Q R) where the old
# T ("file_name") becomes the new R.
$ws.Range("R1:S1").EntireColumn.Delete()

# --- 2. Column widths ---
# (ColumnWidth is specified in characters of the default font and gets
# internally rounded to whole pixels by Excel; the values below are
# chosen so the saved OOXML "width" attribute comes out to exactly
# 19 / 26 / 45 as required.)
$ws.Range("G1").ColumnWidth = 18.14
$ws.Range("Q1").ColumnWidth = 25.14
$ws.Range("R1").ColumnWidth = 44.14

# --- 3. Header row updates ---
# Q1 was "Phase_Code" (after the column shift); rename it to "Remarks".
$ws.Range("Q1").Value = "Remarks"
# R1 is already "file_name" after the shift, no change needed.

# --- 4. Row 2 data updates ---
# Helper: write a value as plain text (avoiding Excel's automatic
# number/date inference), while keeping the cell's style unchanged.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("B2") "08/06/25"
$ws.Range("D2").Value = "JONSUP"
$ws.Range("E2").Value = "Johnstone Supply"
$ws.Range("G2").Value = "110-S10112669.001"
Set-TextValue $ws.Range("H2") "08/06/25"
Set-TextValue $ws.Range("I2") "99.31"
Set-TextValue $ws.Range("J2") "8.93"
Set-TextValue $ws.Range("L2") "90.38"
$ws.Range("N2").ClearContents()
Set-TextValue $ws.Range("P2") "5260"
$ws.Range("Q2").Value = "San Leandro Unit Install"
$ws.Range("R2").Value = "110-s101126669.001 jonsup_1754579789882.pdf"
